# Rename "Sheet1" -> "Providers"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Providers"

# Update the saved selection on the sheet from O6 to A2:XFD5 (active cell A2)
$ws.Range("A2:XFD5").Select()
